$d = $word.ActiveDocument
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Delete()

$p4b = $d.Paragraphs.Item(4)
try {
    $p4b.Range.ParagraphFormat.Reset()
    Write-Host "Reset ok"
} catch {
    Write-Host "Reset failed: $_"
}
Write-Host "ListType after reset: $($p4b.Range.ListFormat.ListType)"
$p4b.Range.ListFormat.RemoveNumbers()
Write-Host "ListType after removenumbers: $($p4b.Range.ListFormat.ListType)"
$p4b.Style = "Listenabsatz"
Write-Host "Text: [$($p4b.Range.Text)]"
